# Checklist.xlsx: Battle section works now (not parallel yet), Auth
# section started. Fill in the "Points" (column B) for the feature rows
# that are now fully done, and leave the viewport where work continued.
#
# - Trading (row 26) -> 3 / 3 points
# - Contains tracked time (row 51) -> 0.5 / 0.5 points
# - Contains link to GIT (row 52) -> 0.5 / 0.5 points
#
# The Sum-Points formulas in rows 56/58 (and their D-column deltas)
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B26").Value = 3
$ws.Range("B51").Value = 0.5
$ws.Range("B52").Value = 0.5

# Scroll/select where work continued (Auth section around row 14).
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B32").Select()
